# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" note with the new rates ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$oldText = $wsHoja1.Range("A1").Value()
$newText = $oldText.Replace("1000 Bs = 8.7 = 36304.35 pesos", "1000 Bs = 8.51 = 35368.09 pesos").Replace("36304.35 pesos = 8.66 = 953.1 Bs", "35368.09 pesos = 8.49 = 960.35 Bs")
$wsHoja1.Range("A1").Value = $newText

# --- tasas: update the "tasas" block (N10/O10/N12/O12) ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 117.5
$wsTasas.Range("O10").Value = 4155.75
$wsTasas.Range("N12").Value = 4166
$wsTasas.Range("O12").Value = 113.12
